# Updating Excel Modification Code + Cucumber Features
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Set 1")
$ws2 = $wb.Worksheets.Item("Set 2")

# --- Set 1 sheet: rename first test case, flip some trigger flags ---
$ws1.Range("B2").Value = "TC1"
$ws1.Range("D2").Value = $false
$ws1.Range("D4").Value = $false
$ws1.Range("D5").Value = $false

# --- Set 2 sheet: flip some trigger flags ---
$ws2.Range("D3").Value = $false
$ws2.Range("D6").Value = $false

# --- Selections / active sheet bookkeeping ---
[void]$ws2.Range("D2:D6").Select()
[void]$ws1.Activate()
[void]$ws1.Range("B2").Select()
